$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the "Results" heading paragraph.
# ------------------------------------------------------------------
$resultsIndex = $null
$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -match "^Results\r?$") {
        $resultsIndex = $i
        break
    }
}

if ($resultsIndex -eq $null) {
    throw "Could not locate the 'Results' heading paragraph."
}

# ------------------------------------------------------------------
# 1) Remove three of the plain empty "spacer" paragraphs (spacing
#    after=0, default line spacing) that sit right before the empty
#    paragraph that immediately precedes the "Results" heading (that
#    last empty paragraph carries explicit line=240/auto and stays).
#    Deleting from the highest index down keeps earlier indices valid.
# ------------------------------------------------------------------
$lineParaIndex = $resultsIndex - 1
for ($k = 0; $k -lt 3; $k++) {
    $d.Paragraphs.Item($lineParaIndex - 1).Range.Delete()
    $resultsIndex = $resultsIndex - 1
    $lineParaIndex = $lineParaIndex - 1
}

# ------------------------------------------------------------------
# 2) Stamp the "Results" run with a <w:lastRenderedPageBreak/> marker
#    immediately before its text, as in the target revision.
#    InsertXML replaces the contents of the range it is called on, so
#    it is invoked on a range spanning exactly the "Results" text,
#    re-supplying that same text together with the new break marker.
# ------------------------------------------------------------------
$resultsRange = $d.Paragraphs.Item($resultsIndex).Range
$textRange = $d.Range($resultsRange.Start, $resultsRange.Start + 7)

if ($textRange.Text -ne "Results") {
    throw "Unexpected range text before InsertXML: [$($textRange.Text)]"
}

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body><w:p><w:r w:rsidRPr="00861F0F"><w:rPr><w:b/><w:bCs/><w:sz w:val="26"/></w:rPr>' +
       '<w:lastRenderedPageBreak/><w:t>Results</w:t></w:r></w:p></w:body></w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

$textRange.InsertXML($xml)
